# Auto commit at 2026-01-17 8:52:41.95
# Refresh the monthly metrics on the "Metrics" sheet with the latest figures.
# Everything that depends on these cells (the "today" sheet's B11:B22 formulas
# that reference Metrics!B2:B13, plus the E/F helper columns, and A1's
# TODAY()-1 snapshot) recalculates automatically once the script finishes.

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")

# Remember which sheet is active so we can restore it after touching the
# selection on "Metrics" (selecting a range requires that sheet to be active).
$originalActive = $wb.ActiveSheet

$metrics.Range("B2").Value = 299765.67
$metrics.Range("B3").Value = 221449.88000000003
$metrics.Range("B4").Value = 77829.60000000002
$metrics.Range("B5").Value = 12299
$metrics.Range("B6").Value = 5935636.3999999985
$metrics.Range("B7").Value = 4992167.5100000007
$metrics.Range("B8").Value = 1741921.42
$metrics.Range("B9").Value = 232576
$metrics.Range("B10").Value = 34401017.389999993
$metrics.Range("B11").Value = 32267442.670000002
$metrics.Range("B12").Value = 12023643.460000001
$metrics.Range("B13").Value = 1330206

# Update the remembered selection on the "Metrics" sheet (it moved from D19
# to E21) without changing which sheet tab is active.
$metrics.Select()
$metrics.Range("E21").Select()
$originalActive.Select()
